$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New full script (Outlook + progressbar) for "Cobranca de materiais WMO para GTA" ---
$script167 = @'
# Default model for SAP automations, developed by Robert Aron Zimmermann, using Google AI Studio tuned prompt model
from sap_functions import SAP
from outlook import Outlook
import progressbar

default_language = 'PT'
login = open('sap_login.txt', 'r').readline().strip().split(',')
scheduled_execution = {'scheduled?': False, 'username': login[0], 'password': login[1], 'principal': '100'}
sap_window = 0


# Cobrança de materiais WMO para GTA
# Automatização desenvolvida para efetuar a cobrança de matérias brutos utilizados na fabricação de peças usinadas
# para alternadores, dos metalúrgicos e usinagem de eixos da WMO.

# Solicitado por Beatriz Silva de Andrade Graciosa
# Desenvolvido por Robert Aron Zimmermann

def read_file(filename):
    clients = []
    with open(filename, 'r') as file:
        lines = file.readlines()
        for line in lines:
            data = line.strip().split(',')
            clients.append({'mrp': data[0], 'name': data[1], 'email': data[2]})
    return clients


def get_copy(filename):
    with open(filename, 'r') as file:
        lines = file.readlines()
        for line in lines:
            return line


class Work:
    def __init__(self):
        self.sap = SAP(sap_window, scheduled_execution, default_language)
        self.orders = ''
        self.materials = []
        self.results = []
        self.clients = read_file('clients.txt')
        self.copy = get_copy('copy_email.txt')

    def COOIS(self):
        self.sap.select_transaction('COOIS')
        self.sap.insert_variant('MATERIAIS GTA')
        self.sap.run_actual_transaction()
        my_grid = self.sap.get_my_grid()
        row_number = self.sap.get_my_grid_count_rows(my_grid)
        for i in range(row_number):
            self.orders = f"{self.orders}{my_grid.getCellValue(i, 'AUFNR')}\n"

    def CO24(self):
        self.sap.select_transaction('CO24')
        self.sap.write_text_field('Layout', '/LISTA_MAT')
        self.sap.multiple_selection_field('Ordem')
        self.sap.multiple_selection_paste_data(self.orders)
        self.sap.run_actual_transaction()
        my_grid = self.sap.get_my_grid()
        row_number = self.sap.get_my_grid_count_rows(my_grid)
        for i in range(row_number):
            if my_grid.getCellValue(i, 'MATNR') not in self.materials:
                self.materials.append(my_grid.getCellValue(i, 'MATNR'))

    def MD04(self, material):
        self.sap.select_transaction('MD04')
        self.sap.write_text_field('Material', material)
        self.sap.write_text_field('Centro', '1200')
        self.sap.run_actual_transaction()
        my_table = self.sap.get_my_table()
        target_column = my_table.getCell(-1, 7).Text

        if target_column == 'Centro fornec./recept.':
            for i in range(my_table.VisibleRowCount):
                center = my_table.getCell(i, 7).Text
                if center != '':
                    self.sap.write_text_field('Centro', str(center))
                    self.sap.run_actual_transaction()

                    my_table = self.sap.get_my_table()
                    for i2 in range(my_table.VisibleRowCount):
                        qty_available = my_table.getCell(0, 5).Text
                        if int(qty_available) > 0:
                            self.sap.press_button('Expandir detalh.cabeç.')
                            plan_mrp = self.sap.get_text_at_side('Planejador MRP', 1)
                            return self.results.append(
                                {'material': material, 'quantity': qty_available, 'center': center,
                                 'mrp': str(plan_mrp)})
                        else:
                            return


if __name__ == '__main__':
    work = Work()
    email = Outlook()
    work.COOIS()
    work.CO24()
    bar = progressbar.ProgressBar(max_value=len(work.materials) - 1)

    for index in range(len(work.materials)):
        work.MD04(str(work.materials[index]))
        bar.update(index)
    texto_email = {}
    for item in work.clients:
        texto_email[item['name']] = ''
    for result in work.results:
        nome_area = next((item['name'] for item in work.clients if item['mrp'] == result['mrp']), None)
        email_area = next((item['email'] for item in work.clients if item['mrp'] == result['mrp']), None)
        if nome_area is not None and email_area is not None:
            texto_email[
                nome_area] = f"{texto_email[nome_area]}{result['material']} possui disponível {result['quantity']} no centro {result['center']}<br>"
        else:
            print(f"O planejador {result['mrp']} não está sendo cobrado! Insira o mesmo no arquivo de texto!")
    body_mail = ''
    receptores = ''
    for item in work.clients:
        if texto_email[item['name']] != '':
            receptores = f"{receptores}{item['email']};"
            body_mail = f"{body_mail}<h4>{item['name']}</h4><h5>{texto_email[item['name']]}</h5><br>"
    email.send_email('Materiais WEN', receptores,
                     f"Bom dia,<br><br>Gentileza enviar pedido abaixo com o lote completo, materiais já se "
                     f"encontram em estoque<br><br>{body_mail}",
                     work.copy)
'@

# --- New description text for the same automation (with the email-sending steps) ---
$desc168 = @'
Título: Cobrança de materiais WMO para GTA
Descrição: Automatização desenvolvida para efetuar a cobrança de matérias brutos utilizados na fabricação de peças usinadas para alternadores, dos metalúrgicos e usinagem de eixos da WMO.
Solicitado por: Beatriz Silva de Andrade Graciosa
Desenvolvido por: Robert Aron Zimmermann

Acessar a transação COOIS
Inserir variante "MATERIAIS GTA"
executar
percorrer tabela inserindo os dados da coluna "AUFNR" em uma variável de texto chamada orders
Criar um dicionário chamado "clients", adicionar dados de um arquivo chamado "clients.txt" separando por virgulas, nesse arquivo existe os campos "mrp","name" e "email"
Criar uma string chamada "copy", percorrer todas as lnhas de um arquivo chamado "copy_email.txt" e ler todas as linhas do mesmo

Acessar transação CO24
no campo "Layout" escrever "/LISTA_MAT"
colar a variável "orders" na seleção múltipla do campo "Ordem"
executar
Percorrer as linhas da tabela verificando se o item do campo "MATNR" não está na lista materials
Se não estiver dentro da lista então acrescentar na mesma

Para cada item de "work.materials":
Acessar transação MD04
inserir no campo "Material" o respectivo material
Inserir no campo "Centro" o texto "1200"
executar transação
Caso o título da coluna 7 do campo Flex seja igual a "Centro fornec./recept.":
Percorrer as linhas visíveis do campo Flex:
armazenar o valor da célula em uma variável chamada "center"
Caso "center" seja diferente de vazio então:
Escrever no campo "Centro" o valor de "center" convertido em texto
executar
Armazenar em uma variável chamada qty_avaiable o valor da linha 0 da coluna 5 do campo Flex convertido em inteiro
Caso qty_avaiable seja maior que 0:
Pressionar o botão "Expandir detalh.cabeç."
armazenar em uma variável chamada "plan_mrp" o texto ao lado de "Planejador MRP"
Armazenar em um dicionário o respectivo material, a variável "qty_avaiable", a variável "center" e "plan_mrp"
Caso não seja maior que 0 então retornar

Para cada item do dicionário criado anteriormente realizar o procedimento abaixo:
coletar o nome da área e email da área fazendo relação do campo "mrp" entre os dicionários "clients" e "results"
caso encontre os dois campos então adicionar o texto em um dicionário chamado "texto_email": "{nome da área}{número material} possui disponível {quantidade} no centro {número do centro}"
Caso não encontre então fazer um print avisando que o respectivo "mrp" não está sendo cobrado
Enviar um email com todos os "email" de "clients", em cópia adicionar a string "copy", o título dele deve ser "Materiais WEN", o corpo do email deve seguir o modelo:
Bom dia,
Gentileza enviar abaixo com o lote completo, materiais já se encontram em estoque
{Percorrer todo o dicionário "texto_email" e  escrever:} {nome do cliente}{texto_email}
'@

# --- Updated "Production Order Release" script (drops the `from outlook import Outlook` import) ---
$script169 = @'
# Default model for SAP automations, developed by Robert Aron Zimmermann, using Google AI Studio tuned prompt model
from sap_functions import SAP
from excel import ExcelHandler

default_language = 'EN'
login = open('sap_login.txt', 'r').readline().strip().split(',')
scheduled_execution = {'scheduled?': False, 'username': login[0], 'password': login[1], 'principal': '100'}
sap_window = 0


# Production Order Release 
# To start Production, orders need to be released using COHV transaction  and before that
# Workflow Must be clear. to track the Work Flow Using ZTMM069.

# Solicitado por Selvaganapathy S
# Desenvolvido por Robert Aron Zimmermann

class Work:
    def __init__(self):
        self.sap = SAP(sap_window, scheduled_execution, default_language)
        self.materials = []
        self.users = []
        self.excel = ExcelHandler('Orders.xlsm')

    def COHV(self):
        self.sap.select_transaction('COHV')
        self.sap.insert_variant('SARANYAM')
        self.sap.run_actual_transaction()
        my_grid = self.sap.get_my_grid()
        rows = self.sap.get_my_grid_count_rows(my_grid)

        for i in range(rows):
            self.materials.append(my_grid.getCellValue(i, 'MATNR'))

    def ZTMM069(self):
        self.sap.select_transaction('ZTMM069')
        self.sap.clean_all_fields()
        self.sap.write_text_field('Interval of days', '999')
        self.sap.write_text_field('Plant', '6200')
        self.sap.write_text_field('Layout', '/RAJBLOCK')
        self.sap.multiple_selection_field('Material')
        self.sap.multiple_selection_paste_data('\n'.join(self.materials))
        self.sap.run_actual_transaction()

        my_grid = self.sap.get_my_grid()
        rows = self.sap.get_my_grid_count_rows(my_grid)
        for i in range(rows):
            if f"{my_grid.getCellValue(i, 'USER')}@weg.net" not in self.users:
                self.users.append(f"{my_grid.getCellValue(i, 'USER')}@weg.net")

        self.excel.load_workbook()
        self.excel.select_sheet('Principal')
        self.excel.clean_data(2, self.excel.count_columns(1), 2, self.excel.count_rows(2))
        self.excel.sap_write_my_grid(my_grid, rows, 1, 2)
        self.excel.save_workbook()
        self.excel.close_workbook()


if __name__ == '__main__':
    mail_outlook = Outlook()
    work = Work()
    work.COHV()
    work.ZTMM069()
    mail_outlook.send_email('Production Order Release', ';'.join(work.users),
                            'The Production Order Release spreadsheet follows...', attachments=f'Orders.xlsm')
'@

# Append the new row (87) with the complete "Cobranca de materiais WMO para GTA" automation.
# Write column B (the script) first, then column A (the description) -- matches the order the
# shared strings were introduced when the author originally saved the workbook.
$ws.Range("B87").Value = $script167
$ws.Range("A87").Value = $desc168

$ws.Range("A87:B87").Style = $ws.Range("A86:B86").Style
$ws.Rows.Item(87).RowHeight = 91.5

# Row 78 column B held the OLD "Production Order Release" script (with the stray `from outlook
# import Outlook`); replace it with the corrected script text.
$ws.Range("B78").Value = $script169

$ws.Range("B87").Select()
